# RPA datasets push 2024-06-20
# Insert 3 new IPO listing rows (2024-06-19) at the top of the data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the 3 blank rows away from the bold header row (at rows 5-7) so the
# bold header formatting does not bleed into the new rows, then move the old
# rows 2-4 data down into rows 5-7 where it belongs after the insert.
$ws.Range("A5:Q7").EntireRow.Insert(-4160)
$ws.Range("A2:Q4").Copy($ws.Range("A5:Q7"))

function Set-TextCell($row, $col, $text) {
    # Force text (not date/number) interpretation, then strip the helper
    # number format back off so the cell ends up with the default style.
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

function Set-NumCell($row, $col, $num) {
    $ws.Cells.Item($row, $col).Value = $num
}

function Set-NewRow($row, $A, $B, $C, $D, $E, $F, $G, $H, $I, $J, $K, $L, $M, $N, $O, $P, $Q) {
    Set-TextCell $row 1  $A
    $ws.Cells.Item($row, 2).Value = $B
    $ws.Cells.Item($row, 3).Value = $C
    Set-NumCell  $row 4  $D
    $ws.Cells.Item($row, 5).Value = $E
    Set-NumCell  $row 6  $F
    $ws.Cells.Item($row, 7).Value = $G
    $ws.Cells.Item($row, 8).Value = $H
    $ws.Cells.Item($row, 9).Value = $I
    $ws.Cells.Item($row, 10).Value = $J
    $ws.Cells.Item($row, 11).Value = $K
    $ws.Cells.Item($row, 12).Value = $L
    Set-NumCell  $row 13 $M
    Set-NumCell  $row 14 $N
    Set-TextCell $row 15 $O
    Set-TextCell $row 16 $P
    Set-NumCell  $row 17 $Q
}

Set-NewRow 2 "2024-06-19" "한국제14호스팩"     "코스닥" 80  "한국" 80  "-" "-" "-" "-" "대표" "-" 2000  100 "2024-06-10" "2024-06-13" 3000000
Set-NewRow 3 "2024-06-19" "미래에셋비전스팩5호" "코스닥" 95  "미래" 95  "-" "-" "-" "-" "대표" "-" 2000  100 "2024-06-10" "2024-06-13" 3562500
Set-NewRow 4 "2024-06-19" "씨어스테크놀로지"   "코스닥" 221 "한국" 221 "-" "-" "-" "-" "대표" "-" 17000 100 "2024-06-10" "2024-06-13" 975000
